# Open the workbook / worksheet that is already active in the session.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The uploaded file changes the value stored in C5 from 25 to 20.
$ws.Range("C5").Value = 20

# The author's last on-screen selection (saved into the sheet view) moved
# from D4 to C5 - mirror that by selecting C5 before saving.
$ws.Range("C5").Select()
